# Update line loading results (pl_mw.xlsx / Sheet1) for the 380 kV case.
# Each entry is (row, column index, new value) for columns B,C,D,F,G,H,I,J,K,N,O
# across rows 2-25. Columns A, E, L, M are left unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @(2, 2, 0.3150526710033148),
    @(2, 3, 0.04542379512027139),
    @(2, 4, 0.3620832968265404),
    @(2, 6, 1.087568053135641),
    @(2, 7, 0.480954622119846),
    @(2, 8, 0.636323509858606),
    @(2, 9, 0.4800662808957519),
    @(2, 10, 0.3480059910106377),
    @(2, 11, 0.3252599649995034),
    @(2, 14, 1.270349397365298),
    @(2, 15, 2.195068470340686),
    @(3, 2, 0.27908169868968),
    @(3, 3, 0.03978548771654289),
    @(3, 4, 0.3502186477398084),
    @(3, 6, 1.087513024883592),
    @(3, 7, 0.4829853969292728),
    @(3, 8, 0.6405549271711237),
    @(3, 9, 0.4848189494474369),
    @(3, 10, 0.3364256760962547),
    @(3, 11, 0.285889301473361),
    @(3, 14, 1.281786794964908),
    @(3, 15, 2.208003653438922),
    @(4, 2, 0.25697828389508),
    @(4, 3, 0.0363072422892543),
    @(4, 4, 0.3430982689277471),
    @(4, 6, 1.088066691611246),
    @(4, 7, 0.4845411660890662),
    @(4, 8, 0.6434062914291871),
    @(4, 9, 0.4879818635149498),
    @(4, 10, 0.3295151997565142),
    @(4, 11, 0.2616713618168376),
    @(4, 14, 1.28922434438465),
    @(4, 15, 2.217122890235188),
    @(5, 2, 0.247967296889783),
    @(5, 3, 0.03488580941944974),
    @(5, 4, 0.340238217505231),
    @(5, 6, 1.088440138634908),
    @(5, 7, 0.4852527932240847),
    @(5, 8, 0.6446319868097063),
    @(5, 9, 0.4893323213400524),
    @(5, 10, 0.3267494220991409),
    @(5, 11, 0.2517918984755738),
    @(5, 14, 1.292359601896333),
    @(5, 15, 2.221135070902875),
    @(6, 2, 0.2464708265418096),
    @(6, 3, 0.03464954114367913),
    @(6, 4, 0.3397658235983272),
    @(6, 6, 1.088511080432511),
    @(6, 7, 0.485375646723341),
    @(6, 8, 0.644839364113821),
    @(6, 9, 0.4895602809434294),
    @(6, 10, 0.3262932057108827),
    @(6, 11, 0.2501508088968052),
    @(6, 14, 1.292886515932199),
    @(6, 15, 2.221819170842139),
    @(7, 2, 0.2568567725790842),
    @(7, 3, 0.03628808851900089),
    @(7, 4, 0.343059528758161),
    @(7, 6, 1.088071129380893),
    @(7, 7, 0.4845504490232244),
    @(7, 8, 0.6434225634228881),
    @(7, 9, 0.4879998270712456),
    @(7, 10, 0.3294776958116898),
    @(7, 11, 0.2615381653864972),
    @(7, 14, 1.289266204794739),
    @(7, 15, 2.217175801268908),
    @(8, 2, 0.3026537954577577),
    @(8, 3, 0.04348313631372491),
    @(8, 4, 0.3579583143229286),
    @(8, 6, 1.087427212241892),
    @(8, 7, 0.4815907105766684),
    @(8, 8, 0.6377299735991571),
    @(8, 9, 0.4816542011853642),
    @(8, 10, 0.3439716294528239),
    @(8, 11, 0.3116945144503802),
    @(8, 14, 1.274206924825297),
    @(8, 15, 2.199284270248171),
    @(9, 2, 0.3923024929096073),
    @(9, 3, 0.05746039656466451),
    @(9, 4, 0.3884744108222549),
    @(9, 6, 1.090822539115138),
    @(9, 7, 0.4782391106418302),
    @(9, 8, 0.6285738568349117),
    @(9, 9, 0.4711525514533967),
    @(9, 10, 0.373980499844734),
    @(9, 11, 0.4096760796036563),
    @(9, 14, 1.247965365929936),
    @(9, 15, 2.173536657080092),
    @(10, 2, 0.4580450681761192),
    @(10, 3, 0.06764599406056959),
    @(10, 4, 0.4116812903696427),
    @(10, 6, 1.096154597249395),
    @(10, 7, 0.4772748272761831),
    @(10, 8, 0.6230672504164829),
    @(10, 9, 0.4646210654230494),
    @(10, 10, 0.3969982699750574),
    @(10, 11, 0.4814088177074041),
    @(10, 14, 1.230687046427317),
    @(10, 15, 2.160312604699243),
    @(11, 2, 0.4879215034524123),
    @(11, 3, 0.07226099587605006),
    @(11, 4, 0.4224084889547441),
    @(11, 6, 1.099196101855881),
    @(11, 7, 0.47716215534345),
    @(11, 8, 0.6208265014320489),
    @(11, 9, 0.4619068695604369),
    @(11, 10, 0.407681162895841),
    @(11, 11, 0.5139816020190722),
    @(11, 14, 1.223260423943092),
    @(11, 15, 2.155532960990314),
    @(12, 2, 0.4992300251639961),
    @(12, 3, 0.07400585069241572),
    @(12, 4, 0.4264949361546826),
    @(12, 6, 1.10043634198793),
    @(12, 7, 0.477166410462516),
    @(12, 8, 0.6200159327262611),
    @(12, 9, 0.4609160387399776),
    @(12, 10, 0.4117569843632936),
    @(12, 11, 0.5263070288214067),
    @(12, 14, 1.220510428673776),
    @(12, 15, 2.153900767746222),
    @(13, 2, 0.4967947677076836),
    @(13, 3, 0.07363018854709935),
    @(13, 4, 0.4256137694521556),
    @(13, 6, 1.100165299840597),
    @(13, 7, 0.4771634064862837),
    @(13, 8, 0.620188815841928),
    @(13, 9, 0.4611277870627823),
    @(13, 10, 0.4108778300754352),
    @(13, 11, 0.5236529474368297),
    @(13, 14, 1.221099918485059),
    @(13, 15, 2.154244383917472),
    @(14, 2, 0.4888519672446705),
    @(14, 3, 0.07240460168077334),
    @(14, 4, 0.4227441977408546),
    @(14, 6, 1.0992963641866),
    @(14, 7, 0.4771615647806726),
    @(14, 8, 0.6207590550373396),
    @(14, 9, 0.4618246121476055),
    @(14, 10, 0.4080158733322321),
    @(14, 11, 0.5149958104258587),
    @(14, 14, 1.223032931547873),
    @(14, 15, 2.155395116945556),
    @(15, 2, 0.4839860961839122),
    @(15, 3, 0.0716535340530271),
    @(15, 4, 0.4209896596642864),
    @(15, 6, 1.098775637385444),
    @(15, 7, 0.477166548581863),
    @(15, 8, 0.6211132845311909),
    @(15, 9, 0.4622562535800334),
    @(15, 10, 0.406266805312498),
    @(15, 11, 0.5096918470628964),
    @(15, 14, 1.224225071536278),
    @(15, 15, 2.156123123333657),
    @(16, 2, 0.4560919059717321),
    @(16, 3, 0.06734401340078477),
    @(16, 4, 0.4109836523831234),
    @(16, 6, 1.095968214859738),
    @(16, 7, 0.4772887541948378),
    @(16, 8, 0.6232190020802051),
    @(16, 9, 0.464803617788494),
    @(16, 10, 0.3963043800447537),
    @(16, 11, 0.4792788645410155),
    @(16, 14, 1.231181113067819),
    @(16, 15, 2.160649839113091),
    @(17, 2, 0.4389714968253884),
    @(17, 3, 0.06469546607814891),
    @(17, 4, 0.4048887634508844),
    @(17, 6, 1.094403647937909),
    @(17, 7, 0.4774472484731902),
    @(17, 8, 0.6245784381252761),
    @(17, 9, 0.4664321828355753),
    @(17, 10, 0.390247021665644),
    @(17, 11, 0.4606059096934416),
    @(17, 14, 1.235559406246811),
    @(17, 15, 2.163743422369691),
    @(18, 2, 0.4291214899940883),
    @(18, 3, 0.06317035784374525),
    @(18, 4, 0.4013991831421606),
    @(18, 6, 1.093561729063438),
    @(18, 7, 0.4775690908195642),
    @(18, 8, 0.6253852232348365),
    @(18, 9, 0.4673930802603188),
    @(18, 10, 0.3867829458575045),
    @(18, 11, 0.4498602289116604),
    @(18, 14, 1.238118476292883),
    @(18, 15, 2.165639111269783),
    @(19, 2, 0.4257859871306664),
    @(19, 3, 0.06265368732279342),
    @(19, 4, 0.4002204317816336),
    @(19, 6, 1.093286630133932),
    @(19, 7, 0.4776156125829871),
    @(19, 8, 0.6256626601404349),
    @(19, 9, 0.4677225776881286),
    @(19, 10, 0.3856134981349868),
    @(19, 11, 0.4462210050238014),
    @(19, 14, 1.238991939869383),
    @(19, 15, 2.166300940352158),
    @(20, 2, 0.4407942884571696),
    @(20, 3, 0.06497758855773839),
    @(20, 4, 0.4055359162576337),
    @(20, 6, 1.094564198768225),
    @(20, 7, 0.4774272008264049),
    @(20, 8, 0.6244311498349617),
    @(20, 9, 0.4662563156092538),
    @(20, 10, 0.3908897726862506),
    @(20, 11, 0.4625942511640062),
    @(20, 14, 1.235089107784635),
    @(20, 15, 2.163402064133649),
    @(21, 2, 0.4911851024842804),
    @(21, 3, 0.07276466149325245),
    @(21, 4, 0.4235864028509013),
    @(21, 6, 1.099549190752896),
    @(21, 7, 0.4771608319407079),
    @(21, 8, 0.6205905322290732),
    @(21, 9, 0.461618934056883),
    @(21, 10, 0.4088556735702014),
    @(21, 11, 0.5175388763275066),
    @(21, 14, 1.222463467891828),
    @(21, 15, 2.155052294471943),
    @(22, 2, 0.5240887370082987),
    @(22, 3, 0.07783791540603602),
    @(22, 4, 0.435524928169599),
    @(22, 6, 1.103322910719385),
    @(22, 7, 0.4772602595607509),
    @(22, 8, 0.6183016739072968),
    @(22, 9, 0.4588036715398758),
    @(22, 10, 0.4207748453517297),
    @(22, 11, 0.5533945225193975),
    @(22, 14, 1.214575042121481),
    @(22, 15, 2.150631329224325),
    @(23, 2, 0.5065303927386822),
    @(23, 3, 0.07513172220615161),
    @(23, 4, 0.4291402310219041),
    @(23, 6, 1.101261642140997),
    @(23, 7, 0.4771821517639552),
    @(23, 8, 0.6195030540742152),
    @(23, 9, 0.4602865026134459),
    @(23, 10, 0.414397139903258),
    @(23, 11, 0.5342628553776763),
    @(23, 14, 1.2187520179057),
    @(23, 15, 2.152896075267677),
    @(24, 2, 0.4399702269070929),
    @(24, 3, 0.0648500485249599),
    @(24, 4, 0.4052432934417709),
    @(24, 6, 1.094491434401007),
    @(24, 7, 0.4774361686746147),
    @(24, 8, 0.6244976602866217),
    @(24, 9, 0.4663357485138455),
    @(24, 10, 0.3905991276965182),
    @(24, 11, 0.4616953541700468),
    @(24, 14, 1.235301599105565),
    @(24, 15, 2.163556027289189),
    @(25, 2, 0.368069797886136),
    @(25, 3, 0.05369363063287835),
    @(25, 4, 0.3800804379075942),
    @(25, 6, 1.089405557563062),
    @(25, 7, 0.4788829379785255),
    @(25, 8, 0.6308362888892631),
    @(25, 9, 0.473785598340065),
    @(25, 10, 0.3656921440548189),
    @(25, 11, 0.383212269017065),
    @(25, 14, 1.254712657564994),
    @(25, 15, 2.179502383522859)
)

foreach ($u in $updates) {
    $row = $u[0]
    $col = $u[1]
    $val = $u[2]
    $ws.Cells.Item($row, $col).Value = $val
}
